$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'polyester athletic pants men'
$ws.Cells.Item(2, 1).Value = 'compression for the knee'
$ws.Cells.Item(3, 1).Value = 'volleyball youth'
$ws.Cells.Item(4, 1).Value = 'knee sleeve with padding'
$ws.Cells.Item(5, 1).Value = 'kneeling pad exercise'
$ws.Cells.Item(6, 1).Value = 'baseball items'
$ws.Cells.Item(7, 1).Value = 'boys cold leggings'
$ws.Cells.Item(8, 1).Value = 'patella knee band'
$ws.Cells.Item(9, 1).Value = 'leg sleeves for basketball youth boys'
$ws.Cells.Item(10, 1).Value = 'mesh capri leggings'
$ws.Cells.Item(11, 1).Value = 'mens black compression pants'
$ws.Cells.Item(12, 1).Value = 'best knee pads construction'
$ws.Cells.Item(13, 1).Value = 'paintball pants small'
$ws.Cells.Item(14, 1).Value = 'below the knee'
$ws.Cells.Item(15, 1).Value = 'football pads'
$ws.Cells.Item(16, 1).Value = 'leg sleeves for basketball youth'
$ws.Cells.Item(17, 1).Value = 'knee pads for working'
$ws.Cells.Item(18, 1).Value = 'volleyball shorts longer length'
$ws.Cells.Item(19, 1).Value = 'kneepads construction'
$ws.Cells.Item(20, 1).Value = 'youth softball pants black'
$ws.Cells.Item(21, 1).Value = 'mountain skin pants'
$ws.Cells.Item(22, 1).Value = 'padded snowboarding shorts'
$ws.Cells.Item(23, 1).Value = 'volleyball kneepads'
$ws.Cells.Item(24, 1).Value = 'boys basketball gear'
$ws.Cells.Item(25, 1).Value = 'baseball leg guard'
$ws.Cells.Item(26, 1).Value = 'construction gel knee pads'
$ws.Cells.Item(27, 1).Value = 'girls youth softball pants black'
$ws.Cells.Item(28, 1).Value = 'running knee band'
$ws.Cells.Item(29, 1).Value = 'recovery pants'
$ws.Cells.Item(30, 1).Value = 'knee foam pad'
$ws.Cells.Item(31, 1).Value = 'lacrosse shorts girls'
$ws.Cells.Item(32, 1).Value = 'compression sleeve knee pads'
$ws.Cells.Item(33, 1).Value = 'compression running tights men'
$ws.Cells.Item(34, 1).Value = 'bjj shorts for men'
$ws.Cells.Item(35, 1).Value = 'football shorts for men'
$ws.Cells.Item(36, 1).Value = 'weight basketball'
$ws.Cells.Item(37, 1).Value = 'professional work knee pads'
$ws.Cells.Item(38, 1).Value = 'wrestling shorts'
$ws.Cells.Item(39, 1).Value = 'mesh basketball shorts for men'
$ws.Cells.Item(40, 1).Value = 'catchers gear leg guard'
$ws.Cells.Item(41, 1).Value = 'gel knee pads construction'
$ws.Cells.Item(42, 1).Value = 'acl knee'
$ws.Cells.Item(43, 1).Value = 'knee sleeve wrestling'
$ws.Cells.Item(44, 1).Value = 'spandex shorts men'
$ws.Cells.Item(45, 1).Value = 'sheer test'
$ws.Cells.Item(46, 1).Value = 'black leggings youth'
$ws.Cells.Item(47, 1).Value = 'calf sleeve padded'
$ws.Cells.Item(48, 1).Value = 'compression hip'
$ws.Cells.Item(49, 1).Value = 'knee sleeve for wrestling'
$ws.Cells.Item(50, 1).Value = 'compression knee for men'
$ws.Cells.Item(51, 1).Value = 'work knee pads'
$ws.Cells.Item(52, 1).Value = 'knees pads'
$ws.Cells.Item(53, 1).Value = 'medium youth baseball pants'
$ws.Cells.Item(54, 1).Value = 'knee compression sleeve with knee pad'
$ws.Cells.Item(55, 1).Value = 'knee pads with gel'
$ws.Cells.Item(56, 1).Value = 'football compression shorts youth'
$ws.Cells.Item(57, 1).Value = 'pants mountain'
$ws.Cells.Item(58, 1).Value = 'hex soccer'
$ws.Cells.Item(59, 1).Value = 'knee guards mountain biking'
$ws.Cells.Item(60, 1).Value = 'construction kneeling pad'
$ws.Cells.Item(61, 1).Value = 'leg guards softball'
$ws.Cells.Item(62, 1).Value = 'sort pants men'
$ws.Cells.Item(63, 1).Value = 'knee pad work'
$ws.Cells.Item(64, 1).Value = 'construction work knee pads'
$ws.Cells.Item(65, 1).Value = 'knee sleeves football'
$ws.Cells.Item(66, 1).Value = 'knee sleeve running men'
$ws.Cells.Item(67, 1).Value = 'squat pants men'
$ws.Cells.Item(68, 1).Value = 'down pants'
$ws.Cells.Item(69, 1).Value = 'adult baseball pants black'
$ws.Cells.Item(70, 1).Value = 'long basketball shorts'
$ws.Cells.Item(71, 1).Value = 'knee padding'
$ws.Cells.Item(72, 1).Value = 'shorts pad'
$ws.Cells.Item(73, 1).Value = 'leg sleeves for men football'
$ws.Cells.Item(74, 1).Value = 'black baseball pants youth'
$ws.Cells.Item(75, 1).Value = 'baseball shorts for boys'
$ws.Cells.Item(76, 1).Value = 'construction knee pads gel'
$ws.Cells.Item(77, 1).Value = 'mountain biking knee pads'
$ws.Cells.Item(78, 1).Value = 'mountain biking pads for men'
$ws.Cells.Item(79, 1).Value = 'work kneepads'
$ws.Cells.Item(80, 1).Value = 'construction knee pads for work'
$ws.Cells.Item(81, 1).Value = 'knee pads for works'
$ws.Cells.Item(82, 1).Value = 'compression knee sleeves for weightlifting'
$ws.Cells.Item(83, 1).Value = 'knee pads for work'
$ws.Cells.Item(84, 1).Value = 'compression shorts football'
$ws.Cells.Item(85, 1).Value = 'under shorts for men'
$ws.Cells.Item(86, 1).Value = 'soccer tights'
$ws.Cells.Item(87, 1).Value = 'knee pad floor'
$ws.Cells.Item(88, 1).Value = 'gel work knee pads'
$ws.Cells.Item(89, 1).Value = 'knee pads for kneeling'
$ws.Cells.Item(90, 1).Value = 'knees pads for work'
$ws.Cells.Item(91, 1).Value = 'biking capris'
$ws.Cells.Item(92, 1).Value = 'hip protector'
$ws.Cells.Item(93, 1).Value = 'volleyball kneepads women'
$ws.Cells.Item(94, 1).Value = 'compression shorts bjj'
$ws.Cells.Item(95, 1).Value = 'basketball sleeve youth boys'
$ws.Cells.Item(96, 1).Value = 'yoga tights'
$ws.Cells.Item(97, 1).Value = 'padded sliding shorts women'
$ws.Cells.Item(98, 1).Value = 'tight leggings'
$ws.Cells.Item(99, 1).Value = 'catchers gear women'
$ws.Cells.Item(100, 1).Value = 'softball gear for men'
